$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from an existing header cell (H1) to the new header cells
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @{
    2  = @(8, 8)
    3  = @(7, 8)
    4  = @(8, 8)
    5  = @(6, 6)
    6  = @(6, 7)
    7  = @(5, 6)
    8  = @(9, 9)
    9  = @(9, 9)
    10 = @(6, 6)
    11 = @(8, 8)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
